$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Session 7 data entry: fill in the missing Lab 7 (Branch&Bound) score
$ws.Range("H6").Value = 10

# Add feedback comment for row 7 (Lab 7 column) first so it becomes shared string 23
$ws.Range("H7").Value = "Good one"

# Mark the test as passed (becomes shared string 24)
$ws.Range("I6").Value = "Passsed"

# Update selection to match the author's final cursor position
$ws.Range("I7:I14").Select()
